# Weekly update: insert a new price record at row 14 for "Camote"
# (Vega Modelo de Temuco) and shift the existing rows 14-34 down to
# rows 15-35, matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14; pushes old rows 14..34 down to 15..35
# and keeps the row-14 formatting (date style on column D, etc.).
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new record.
$ws.Range("A14").Value = 10
$ws.Range("B14").Value = "Vega Modelo de Temuco"
$ws.Range("C14").Value = "La Araucanía"
$ws.Range("D14").Value = 44525
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 100114002
$ws.Range("G14").Value = "Camote"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 40
$ws.Range("K14").Value = 20000
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = 20000
$ws.Range("N14").Value = "$/caja 15 kilos granel"
$ws.Range("O14").Value = "Perú"
$ws.Range("P14").Value = 1333
$ws.Range("Q14").Value = 15
$ws.Range("R14").Value = "Hortaliza"
